# Add a new "Mass" column (F) to the "goods" query table on the active sheet,
# populate it with per-row mass values, and leave the selection where Excel
# would after such an edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the existing query table ("goods") by one column so the table range,
# autofilter and tableColumns collection all expand together (A1:E22 -> A1:F22).
$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()

# Header
$ws.Range("F1").Value = "Mass"

# Data rows (row -> mass value)
$ws.Range("F2").Value = 1000
$ws.Range("F3").Value = 1000
$ws.Range("F4").Value = 1000
$ws.Range("F5").Value = 1000
$ws.Range("F6").Value = 1000
$ws.Range("F7").Value = 1000
$ws.Range("F8").Value = 1000
$ws.Range("F9").Value = 1000
$ws.Range("F10").Value = 1000
$ws.Range("F11").Value = 1000
$ws.Range("F12").Value = 500000
$ws.Range("F13").Value = 50000000
$ws.Range("F14").Value = 1000
$ws.Range("F15").Value = 1000
$ws.Range("F16").Value = 1000
$ws.Range("F17").Value = 1000
$ws.Range("F18").Value = 1000
$ws.Range("F19").Value = 2500
$ws.Range("F20").Value = 1000
$ws.Range("F21").Value = 100
$ws.Range("F22").Value = 100

# Match the numeric formatting already used by the other data columns.
$ws.Range("F2:F22").NumberFormat = "General"

# Leave the selection on the cell just below the newly-filled column, matching
# where Excel lands after entering the last value in a freshly added column.
$ws.Range("F23").Select()
